$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title table cell: "Izvještaj - NASP" -> "Izvještaj - RASUS"
# ---------------------------------------------------------------------------
$r1 = $d.Content
$null = $r1.Find.Execute("NASP", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "RASUS", 2)

# ---------------------------------------------------------------------------
# 2) Split the large "Kada neki čvor..." paragraph: insert a new Heading 2
#    paragraph "Izgradnja stabla" right before it, then patch a handful of
#    words/sentences inside the (now following) body paragraph.
# ---------------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Kada neki ")) {
        $targetIndex = $i
        break
    }
}

$targetPara = $d.Paragraphs.Item($targetIndex)
$targetPara.Range.InsertParagraphBefore()

# the freshly-inserted empty paragraph now sits right before the big one
$headingPara = $d.Paragraphs.Item($targetIndex)
$headingPara.Range.Text = "Izgradnja stabla"
$headingPara.Style = "Heading 2"

# re-fetch the body paragraph (index shifted by one after the insert above)
$bodyPara = $d.Paragraphs.Item($targetIndex + 1)
$bodyRange = $bodyPara.Range

$null = $bodyRange.Find.Execute("prošili", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "proširili", 2)

$null = $bodyRange.Find.Execute( `
    "Na taj način čvor roditelj čini listu", $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    "Na temelju tih potvrdnih poruka  čvor roditelj uvrštava svoju djecu u listu", 2)

$null = $bodyRange.Find.Execute( `
    "stvaranju istog stabla. U tom", $true, $false, $false, $false, $false, `
    $true, 1, $false, "stvaranju istog stabla (isti ID). U tom", 2)

$null = $bodyRange.Find.Execute( `
    "obilazak po mreži. ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "obilazak po mreži.", 2)

# ---------------------------------------------------------------------------
# 3) Append new sentence about newcomers joining future trees to the
#    "Kako je mreža dinamična..." paragraph, with the _GoBack bookmark
#    sitting in the middle of "preplavljivanja" (matches the source edit).
# ---------------------------------------------------------------------------
$endRange = $d.Content
$null = $endRange.Find.Execute( `
    "stablo po kojem je poruka poslana može biti obrisano.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endRange.Collapse(0)
$insertStart = $endRange.Start
$endRange.InsertAfter( `
    " Ovo otvara prostor da novi č" + `
    "vor, koji je tek došao u mrežu, možda postane dio nekih stabala koja će biti generiran" + `
    "a u svrhu budućeg preplavljivanja porukama od strane određenih korijenskih čvora.")
$insertEnd = $endRange.End

# find the freshly-inserted "preplav" (scoped, so we don't hit the unrelated
# "preplavljivanje" mention in the "Opis zadatka" section) and drop the
# _GoBack bookmark right in the middle of "preplav|ljivanja".
$scoped = $d.Range($insertStart, $insertEnd)
$null = $scoped.Find.Execute("preplav", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
$scoped.Collapse(0)
$d.Bookmarks.Add("_GoBack", $scoped)

Write-Output "done"
